# edit.ps1
# 1. Fixes the selection/active-tab state left over on "LoanCollectionTransfer" /
#    "BondConfiguration" (the stray xpath/tabSelected artefact).
# 2. Adds a new "CentreGroupRegistration" module sheet (copied from the
#    LoanCollectionTransfer template sheet so styles/column widths/page setup
#    match the rest of the workbook) and fills in its header/value rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Part 1: Add the new "CentreGroupRegistration" worksheet, placed immediately
# after "BondConfiguration" (the last sheet in the workbook).
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item("LoanCollectionTransfer")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$template.Copy([System.Reflection.Missing]::Value, $lastSheet)

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "CentreGroupRegistration"

# Clear any leftover data from the copied template sheet before writing the
# real content for this module.
$newSheet.Cells.Clear()

# Header row
$newSheet.Cells.Item(1, 1).Value = "TestScenario"
$newSheet.Cells.Item(1, 2).Value = "Run"
$newSheet.Cells.Item(1, 3).Value = "pcRegFormName"
$newSheet.Cells.Item(1, 4).Value = "pcRegFormPcName"
$newSheet.Cells.Item(1, 5).Value = "groupName"
$newSheet.Cells.Item(1, 6).Value = "custName"
$newSheet.Cells.Item(1, 7).Value = "custName2"

# Value row
$newSheet.Cells.Item(2, 1).Value = "Centre Group Registration"
$newSheet.Cells.Item(2, 2).Value = "Yes"
$newSheet.Cells.Item(2, 3).Value = "qwerty"
$newSheet.Cells.Item(2, 4).Value = "zxcvb"
$newSheet.Cells.Item(2, 5).Value = "Group"
$newSheet.Cells.Item(2, 6).Value = "q"
$newSheet.Cells.Item(2, 7).Value = "a"

$newSheet.Range("A1:G1").RowHeight = 45
$newSheet.Range("A2:G2").RowHeight = 60

# This is now the active/selected sheet/tab, with the cursor on J11 (matching
# the state captured when the module was saved).
$newSheet.Activate()
$newSheet.Range("J11").Select()

# ---------------------------------------------------------------------------
# Part 2: "Fix the xpath issue" on LoanCollectionTransfer / BondConfiguration
# — BondConfiguration is no longer the active tab and its lingering selection
# (F7) is replaced with the A1:D2 block.
# ---------------------------------------------------------------------------
$bondConfig = $wb.Worksheets.Item("BondConfiguration")
$bondConfig.Activate()
$bondConfig.Range("A1:D2").Select()

# Re-activate the new sheet so it ends up as the saved active tab, matching
# the target workbook state (activeTab points at CentreGroupRegistration).
$newSheet.Activate()
